$wb = $excel.ActiveWorkbook

# Sheet "Impact Matrix": shorten event names in B2 and B3
$wsImpact = $wb.Worksheets.Item("Impact Matrix")
$wsImpact.Range("B2").Value = "NDPS 2026-2030 Launch"
$wsImpact.Range("B3").Value = "IPS / Ethiopay Launch"

# Sheet "Events Metadata": shorten event names in B12 and B13
$wsEvents = $wb.Worksheets.Item("Events Metadata")
$wsEvents.Range("B12").Value = "NDPS 2026-2030 Launch"
$wsEvents.Range("B13").Value = "IPS / Ethiopay Launch"
